$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45186 (2023-09-17) to 45188 (2023-09-19) for every data row (2-181).
for ($r = 2; $r -le 181; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
